$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 621, shifting existing rows 621-670 down to 622-671.
$ws.Rows.Item(621).Insert()

# Populate the newly inserted row 621 with its data (matches prior row's
# unchanged fields, with updated date / volume / price fields).
$ws.Cells.Item(621, 1).Value = 3
$ws.Cells.Item(621, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(621, 3).Value = "Coquimbo"
$ws.Cells.Item(621, 4).Value = 45013
$ws.Cells.Item(621, 5).Value = 5
$ws.Cells.Item(621, 6).Value = 100112037
$ws.Cells.Item(621, 7).Value = "Cebollín"
$ws.Cells.Item(621, 8).Value = "Sin especificar"
$ws.Cells.Item(621, 9).Value = "Primera"
$ws.Cells.Item(621, 10).Value = 230
$ws.Cells.Item(621, 11).Value = 3800
$ws.Cells.Item(621, 12).Value = 4000
$ws.Cells.Item(621, 13).Value = 3904
$ws.Cells.Item(621, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(621, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(621, 16).Value = 108
$ws.Cells.Item(621, 17).Value = 36
$ws.Cells.Item(621, 18).Value = "Hortaliza"
